# Generate Report for Handback
# The handback transform failed for ba93f843-423c-4baa-b122-95177565a56b,
# because the handback file name did not match the handoff file name it
# was supposed to correspond to. Update the Status for that row on the
# Overview roll-up plus both language sheets, and record the Error Detail
# explaining why on each language sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 7 on every sheet corresponds to file ba93f843-423c-4baa-b122-95177565a56b.
# Status moves from "Ready for handoff" to "Handback transform failed" -
# on the Overview sheet this is reflected per-language in columns B (zh-cn)
# and C (de-de), and on each language sheet in its own Status column C.
$overview.Range("B7").Value = "Handback transform failed"
$overview.Range("C7").Value = "Handback transform failed"
$zhcn.Range("C7").Value = "Handback transform failed"
$dede.Range("C7").Value = "Handback transform failed"

# Error Detail (column L) explains the handback/handoff file name mismatch,
# one message per locale referencing that locale's handoff file.
$zhcn.Range("L7").Value = "Handback file name: wkjgjf3d.s02 is different with handoff file name: ba93f843-423c-4baa-b122-95177565a56b.433ef8f43149432a8d7136209652e727894f47f8.zh-cn."
$dede.Range("L7").Value = "Handback file name: wkjgjf3d.s02 is different with handoff file name: ba93f843-423c-4baa-b122-95177565a56b.433ef8f43149432a8d7136209652e727894f47f8.de-de."
